# repull data, push all data, mean calculation
# Update column F (dSF) values on Sheet1 to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -5
    4  = -4
    6  = -3
    7  = -2
    8  = 1
    11 = 7
    16 = 1
    17 = -1
    18 = 1
    19 = -6
    20 = -1
    21 = 3
    22 = -3
    23 = 1
    24 = -2
    26 = 5
    27 = 2
    28 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
